# Generate Report for Archive
#
# The "Status" column value "Ready for handoff" moves to "In Translation"
# on all three sheets (Overview columns E/F hold the per-language status,
# zh-cn/de-de column C holds the same status), and the now-narrower text
# shrinks that column's width accordingly.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns are E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Columns("E:F").ColumnWidth = 12.5

# --- zh-cn sheet: status is column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Columns("C:C").ColumnWidth = 12.5

# --- de-de sheet: status is column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Columns("C:C").ColumnWidth = 12.5
